# edit.ps1 - applies the "Test connection to DB and Display to CSHTML with
# Controller" revision to Plan.docx:
#   1. Trims the run of ~21 blank paragraphs above the "Create Project"
#      heading down to 3, and drops the stale <w:lastRenderedPageBreak/>
#      hint on that heading's run.
#   2. Fills in the trailing blank "ConnectionString" bullet under
#      "Create Project" with the DB-connection test notes, adds a new
#      "Front-end Design:" heading with an "HTML / CSS" bullet, and
#      restores the trailing blank paragraphs after it.

$d = $word.ActiveDocument

$wNs       = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
             '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
             '<pkg:part pkg:name="/part" pkg:contentType="application/xml"><pkg:xmlData>'
$xmlFooter = '</pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# Step 1: collapse the 21 blank paragraphs above "Create Project" to 3.
# ---------------------------------------------------------------------
$createProjPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text.TrimEnd() -eq "Create Project:") {
        $createProjPara = $pp
        break
    }
}

$blankCount = 0
$scanIndex = $createProjPara.Index - 1
while ($scanIndex -ge 1) {
    $pp = $d.Paragraphs.Item($scanIndex)
    if ($pp.Range.Text.TrimEnd() -ne "") { break }
    $blankCount++
    $scanIndex--
}
$keepBlank = 4
$toRemove = $blankCount - $keepBlank
if ($toRemove -gt 0) {
    $firstRemove = $d.Paragraphs.Item($scanIndex + 1)
    $lastRemove  = $d.Paragraphs.Item($scanIndex + $toRemove)
    $delRange = $d.Range($firstRemove.Range.Start, $lastRemove.Range.End)
    $delRange.Delete()
}

# ---------------------------------------------------------------------
# Step 2: rewrite the "Create Project" heading without
# <w:lastRenderedPageBreak/>, keeping its two runs ("Create Project" +
# ":") separate.
# ---------------------------------------------------------------------
$createProjPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text.TrimEnd() -eq "Create Project:") {
        $createProjPara = $pp
        break
    }
}
$headingXml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr>" +
              "<w:r><w:t>Create Project</w:t></w:r>" +
              "<w:r><w:t>:</w:t></w:r></w:p>"
$createProjPara.Range.InsertXML($xmlHeader + $headingXml + $xmlFooter)

# ---------------------------------------------------------------------
# Step 3: the trailing blank "ListParagraph" bullet (ilvl 0 / numId 1)
# right after "Build Model -> Entity Framework" gets replaced by the new
# "ConnectionString" sub-bullets, a new "Front-end Design:" section, and
# the populated "HTML / CSS" bullet - followed by the same run of blank
# paragraphs it used to be followed by.
# ---------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Style.NameLocal -eq "List Paragraph" -and $pp.Range.Text.TrimEnd() -eq "") {
        $targetPara = $pp
    }
}

function ListItemXml([string]$ilvl, [string]$text) {
    return "<w:p $wNs><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr>" +
           "<w:ilvl w:val=`"$ilvl`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr>" +
           "<w:r><w:t>$text</w:t></w:r></w:p>"
}

$emptyP = "<w:p $wNs/>"

$connectionString = "ConnectionString:"
$giangNote        = "Giang: cho ra appsettings.json"
$datNote          = [string]::Join(":", @([char]0x0110 + [char]0x1EA1 + "t", " s" + [char]0x1EED + "a gitignore"))
$datNote          = "$([char]0x0110)$([char]0x1EA1)t: s$([char]0x1EED)a gitignore"
$runNote          = "ch$([char]0x1EA1)y th$([char]0x1EED) xem l$([char]0x1ED7)i g$([char]0x00EC) ko"

$newBullets  = ListItemXml "1" $connectionString
$newBullets += ListItemXml "2" $giangNote
$newBullets += ListItemXml "2" $datNote
$newBullets += ListItemXml "3" $runNote
$newBullets += $emptyP + $emptyP + $emptyP
$newBullets += "<w:p $wNs><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr>" +
               "<w:r><w:t>Front-end Design:</w:t></w:r></w:p>"

$layTu = "L$([char]0x1EA5)y t$([char]0x1EEB) "

$htmlCssXml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr>" +
              "<w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr>" +
              "<w:r><w:rPr><w:b/><w:bCs/><w:color w:val=`"0070C0`"/></w:rPr><w:t>HTML</w:t></w:r>" +
              "<w:r><w:rPr><w:b/><w:bCs/><w:color w:val=`"0070C0`"/></w:rPr><w:t xml:space=`"preserve`"> / CSS</w:t></w:r>" +
              "<w:r><w:t xml:space=`"preserve`">: $layTu</w:t></w:r>" +
              "<w:r><w:rPr><w:b/><w:bCs/><w:color w:val=`"00B0F0`"/></w:rPr><w:t>Project Java</w:t></w:r>" +
              "</w:p>"

# Four blank paragraphs should remain after the "HTML / CSS" bullet (and
# before the pre-existing trailing blank paragraph). InsertXML silently
# folds the very last paragraph mark of the inserted block into whatever
# follows when that last paragraph is itself empty, so one extra blank
# paragraph is appended here purely to compensate for that loss.
$trailingBlanks = $emptyP + $emptyP + $emptyP + $emptyP + $emptyP

$blockXml = $newBullets + $htmlCssXml + $trailingBlanks
$targetPara.Range.InsertXML($xmlHeader + $blockXml + $xmlFooter)
